$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.096772333333333
$ws.Range("H2").Value = 3.290317
$ws.Range("I2").Value = 0.2426185621302128
$ws.Range("J2").Value = 0.2426185621302128
$ws.Range("M2").Value = 8.308763666666666
$ws.Range("N2").Value = 24.926291
$ws.Range("O2").Value = 0.3217360040147476
$ws.Range("P2").Value = 0.3217360040147476
$ws.Range("Q2").Value = 9.112822113805221
$ws.Range("R2").Value = 82.01539902424699
$ws.Range("S2").Value = 0.07805912667957843
$ws.Range("T2").Value = 0.07805912667957843

# Row 3
$ws.Range("G3").Value = 1.096772333333333
$ws.Range("H3").Value = 3.290317
$ws.Range("I3").Value = 0.2426185621302128
$ws.Range("J3").Value = 0.2426185621302128
$ws.Range("M3").Value = 8.621912666666665
$ws.Range("O3").Value = 0.3338619125088609
$ws.Range("P3").Value = 0.3338619125088609
$ws.Range("Q3").Value = 9.456275273216219
$ws.Range("R3").Value = 85.10647745894599
$ws.Range("S3").Value = 0.08100109716294275
$ws.Range("T3").Value = 0.08100109716294275

# Row 4
$ws.Range("G4").Value = 1.096772333333333
$ws.Range("H4").Value = 3.290317
$ws.Range("I4").Value = 0.2426185621302128
$ws.Range("J4").Value = 0.2426185621302128
$ws.Range("M4").Value = 5.000319666666667
$ws.Range("N4").Value = 15.000959
$ws.Range("O4").Value = 0.1936248198758919
$ws.Range("P4").Value = 0.1936248198758919
$ws.Range("Q4").Value = 5.484212268222556
$ws.Range("R4").Value = 49.357910414003
$ws.Range("S4").Value = 0.04697697539101034
$ws.Range("T4").Value = 0.04697697539101034

# Row 5
$ws.Range("G5").Value = 1.096772333333333
$ws.Range("H5").Value = 3.290317
$ws.Range("I5").Value = 0.2426185621302128
$ws.Range("J5").Value = 0.2426185621302128
$ws.Range("M5").Value = 3.893790666666666
$ws.Range("N5").Value = 11.681372
$ws.Range("O5").Value = 0.1507772636004996
$ws.Range("P5").Value = 0.1507772636004996
$ws.Range("Q5").Value = 4.270601874991555
$ws.Range("R5").Value = 38.43541687492399
$ws.Range("S5").Value = 0.03658136289668128
$ws.Range("T5").Value = 0.03658136289668128

# Row 6
$ws.Range("I6").Value = 0.03766810132102297
$ws.Range("J6").Value = 0.03766810132102297
$ws.Range("M6").Value = 8.308763666666666
$ws.Range("N6").Value = 24.926291
$ws.Range("O6").Value = 0.3217360040147476
$ws.Range("P6").Value = 0.3217360040147476
$ws.Range("Q6").Value = 1.414824585923667
$ws.Range("R6").Value = 12.733421273313
$ws.Range("S6").Value = 0.01211918439784856
$ws.Range("T6").Value = 0.01211918439784856

# Row 7
$ws.Range("I7").Value = 0.03766810132102297
$ws.Range("J7").Value = 0.03766810132102297
$ws.Range("M7").Value = 8.621912666666665
$ws.Range("O7").Value = 0.3338619125088609
$ws.Range("P7").Value = 0.3338619125088609
$ws.Range("Q7").Value = 1.468147910792666
$ws.Range("S7").Value = 0.01257594434761428
$ws.Range("T7").Value = 0.01257594434761428

# Row 8
$ws.Range("I8").Value = 0.03766810132102297
$ws.Range("J8").Value = 0.03766810132102297
$ws.Range("M8").Value = 5.000319666666667
$ws.Range("N8").Value = 15.000959
$ws.Range("O8").Value = 0.1936248198758919
$ws.Range("P8").Value = 0.1936248198758919
$ws.Range("Q8").Value = 0.8514594331596668
$ws.Range("R8").Value = 7.663134898437002
$ws.Range("S8").Value = 0.007293479333349917
$ws.Range("T8").Value = 0.007293479333349917

# Row 9
$ws.Range("I9").Value = 0.03766810132102297
$ws.Range("J9").Value = 0.03766810132102297
$ws.Range("M9").Value = 3.893790666666666
$ws.Range("N9").Value = 11.681372
$ws.Range("O9").Value = 0.1507772636004996
$ws.Range("P9").Value = 0.1507772636004996
$ws.Range("Q9").Value = 0.6630385685106667
$ws.Range("R9").Value = 5.967347116596001
$ws.Range("S9").Value = 0.005679493242210206
$ws.Range("T9").Value = 0.005679493242210206

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5018676666666667
$ws.Range("H10").Value = 1.505603
$ws.Range("I10").Value = 0.1110188577571507
$ws.Range("J10").Value = 0.1110188577571507
$ws.Range("M10").Value = 8.308763666666666
$ws.Range("N10").Value = 24.926291
$ws.Range("O10").Value = 0.3217360040147476
$ws.Range("P10").Value = 0.3217360040147476
$ws.Range("Q10").Value = 4.169899834274777
$ws.Range("R10").Value = 37.529098508473
$ws.Range("S10").Value = 0.03571876366506734
$ws.Range("T10").Value = 0.03571876366506733

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5018676666666667
$ws.Range("H11").Value = 1.505603
$ws.Range("I11").Value = 0.1110188577571507
$ws.Range("J11").Value = 0.1110188577571507
$ws.Range("M11").Value = 8.621912666666665
$ws.Range("O11").Value = 0.3338619125088609
$ws.Range("P11").Value = 0.3338619125088609
$ws.Range("Q11").Value = 4.327059192223778
$ws.Range("R11").Value = 38.943532730014
$ws.Range("S11").Value = 0.03706496817535154
$ws.Range("T11").Value = 0.03706496817535152

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5018676666666667
$ws.Range("H12").Value = 1.505603
$ws.Range("I12").Value = 0.1110188577571507
$ws.Range("J12").Value = 0.1110188577571507
$ws.Range("M12").Value = 5.000319666666667
$ws.Range("N12").Value = 15.000959
$ws.Range("O12").Value = 0.1936248198758919
$ws.Range("P12").Value = 0.1936248198758919
$ws.Range("Q12").Value = 2.509498763697445
$ws.Range("R12").Value = 22.585488873277
$ws.Range("S12").Value = 0.02149600633605557
$ws.Range("T12").Value = 0.02149600633605556

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5018676666666667
$ws.Range("H13").Value = 1.505603
$ws.Range("I13").Value = 0.1110188577571507
$ws.Range("J13").Value = 0.1110188577571507
$ws.Range("M13").Value = 3.893790666666666
$ws.Range("N13").Value = 11.681372
$ws.Range("O13").Value = 0.1507772636004996
$ws.Range("P13").Value = 0.1507772636004996
$ws.Range("Q13").Value = 1.954167636368445
$ws.Range("R13").Value = 17.587508727316
$ws.Range("S13").Value = 0.01673911958067628
$ws.Range("T13").Value = 0.01673911958067628

# Row 14
$ws.Range("G14").Value = 2.751641333333334
$ws.Range("H14").Value = 8.254924000000001
$ws.Range("I14").Value = 0.6086944787916135
$ws.Range("J14").Value = 0.6086944787916135
$ws.Range("M14").Value = 8.308763666666666
$ws.Range("N14").Value = 24.926291
$ws.Range("O14").Value = 0.3217360040147476
$ws.Range("P14").Value = 0.3217360040147476
$ws.Range("Q14").Value = 22.86273753409822
$ws.Range("R14").Value = 205.764637806884
$ws.Range("S14").Value = 0.1958389292722532
$ws.Range("T14").Value = 0.1958389292722532

# Row 15
$ws.Range("G15").Value = 2.751641333333334
$ws.Range("H15").Value = 8.254924000000001
$ws.Range("I15").Value = 0.6086944787916135
$ws.Range("J15").Value = 0.6086944787916135
$ws.Range("M15").Value = 8.621912666666665
$ws.Range("O15").Value = 0.3338619125088609
$ws.Range("P15").Value = 0.3338619125088609
$ws.Range("Q15").Value = 23.72441126599022
$ws.Range("R15").Value = 213.519701393912
$ws.Range("S15").Value = 0.2032199028229524
$ws.Range("T15").Value = 0.2032199028229524

# Row 16
$ws.Range("G16").Value = 2.751641333333334
$ws.Range("H16").Value = 8.254924000000001
$ws.Range("I16").Value = 0.6086944787916135
$ws.Range("J16").Value = 0.6086944787916135
$ws.Range("M16").Value = 5.000319666666667
$ws.Range("N16").Value = 15.000959
$ws.Range("O16").Value = 0.1936248198758919
$ws.Range("P16").Value = 0.1936248198758919
$ws.Range("Q16").Value = 13.75908627467956
$ws.Range("R16").Value = 123.831776472116
$ws.Range("S16").Value = 0.117858358815476
$ws.Range("T16").Value = 0.117858358815476

# Row 17
$ws.Range("G17").Value = 2.751641333333334
$ws.Range("H17").Value = 8.254924000000001
$ws.Range("I17").Value = 0.6086944787916135
$ws.Range("J17").Value = 0.6086944787916135
$ws.Range("M17").Value = 3.893790666666666
$ws.Range("N17").Value = 11.681372
$ws.Range("O17").Value = 0.1507772636004996
$ws.Range("P17").Value = 0.1507772636004996
$ws.Range("Q17").Value = 10.71431534174756
$ws.Range("R17").Value = 96.42883807572801
$ws.Range("S17").Value = 0.0917772878809318
$ws.Range("T17").Value = 0.0917772878809318

